# "added duration and moment.js"
#
# 1) Remove the _GoBack bookmark from the "56.48 x 17 ... = 960.16 minutes
#    between 3 and 7" paragraph (it moves to the very end of the doc).
# 2) Split the "960.16/ 60" run into three runs: "960.16", " ", "/ 60".
# 3) Append a blank paragraph, an explanation paragraph (with spell-check
#    proofErr wrapping on the camelCase identifiers), another blank
#    paragraph, and a "Moment.duration().asMinutes()" paragraph that now
#    carries the _GoBack bookmark.

$d = $word.ActiveDocument

function Set-ParagraphBodyXml($paragraph, [string]$innerXml) {
    $r = $paragraph.Range
    $full = $d.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full.InsertXML($pkg)
}

# --- Step 1: find the two paragraphs we need to touch by their text ---
$calcPara = $null
$minutesPara = $null
$lastPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "56.48 x 17*") { $calcPara = $p }
    if ($p.Range.Text -like "960.16/ 60*") { $minutesPara = $p }
}
$lastPara = $d.Paragraphs.Last

# --- Step 1: drop the bookmark from the "56.48 x 17 ..." paragraph ---
Set-ParagraphBodyXml $calcPara (
    '<w:r><w:t>56.48 x 17</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>=</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> 960.16 minutes between 3 and 7</w:t></w:r>'
)

# --- Step 2: split "960.16/ 60" into three runs ---
Set-ParagraphBodyXml $minutesPara (
    '<w:r><w:t>960.16</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>/ 60</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> = 16 hours meaning that the train arrives at 7pm.</w:t></w:r>'
)

# --- Step 3: append the new paragraphs after the last paragraph ---
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newBodyXml =
    '<w:p><w:pPr/></w:p>' +
    '<w:p>' +
        '<w:r><w:t xml:space="preserve">For </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>minutesAway</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> subtract </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>arrivalTime</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">from </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>startTime</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> divided by frequency</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> in minutes</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">. The remainder should equal </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>minutesAway</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr/></w:p>' +
    '<w:p>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Moment.duration</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>().</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/><w:r><w:t>asMinutes</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>()</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($pkg)

Write-Output "edit complete"
